$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pre-Alert Template Import")

$ws.Range("A3").Value = "'JSSO1000252"
$ws.Range("B3").Value = "'JSSO1000252"
$ws.Range("C3").Value = "'JSSO1000252"
$ws.Range("AJ3").Value = "JSCN1000252"
$ws.Range("AN3").Value = "'MBLJSSO1000252"
$ws.Range("AO3").Value = "'HBLJSSO1000252"
